# Applies the crypto price/volume update described in the commit diff.
# Source data cells are stored as text (t="inlineStr"/shared string), even
# when the text looks like a plain number (e.g. "209.26"), so for any new
# value that Excel would otherwise auto-convert to a numeric type we prefix
# it with a leading apostrophe, which forces Excel to keep it as text while
# entering it through the object model, exactly like typing it in the UI.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Bitcoin
$ws.Range('D2').Value = '25.868.73'
$ws.Range('E2').Value = '  -0.46%  '
# Ethereum
$ws.Range('D3').Value = '1.599.12'
$ws.Range('E3').Value = '  -2.09%  '
# TetherUSD
$ws.Range('E4').Value = '  +0.08%  '
# BNB
$ws.Range('D5').Value = '''209.26'
$ws.Range('E5').Value = '  -2.22%  '
# USDC
$ws.Range('E6').Value = '  +0.09%  '
# XRP
$ws.Range('D7').Value = '''0.478'
$ws.Range('E7').Value = '  -5.28%  '
# Cardano
$ws.Range('D8').Value = '''0.245'
$ws.Range('E8').Value = '  -2.79%  '
# Dogecoin
$ws.Range('E9').Value = '  -2.07%  '
# Solana
$ws.Range('D10').Value = '''17.87'
$ws.Range('E10').Value = '  -3.65%  '
# TRON
$ws.Range('E11').Value = '  -0.43%  '
# WrappedliquidstakedEther2.0
$ws.Range('D12').Value = '1.820.26'
$ws.Range('E12').Value = '  -2.14%  '
# WrappedEther
$ws.Range('D13').Value = '1.602.01'
$ws.Range('E13').Value = '  -2.87%  '
# Polkadot
$ws.Range('E14').Value = '  -3.56%  '
# Polygon
$ws.Range('D15').Value = '''0.509'
$ws.Range('E15').Value = '  -3.90%  '
# WrappedBTC
$ws.Range('D16').Value = '25.854.49'
$ws.Range('E16').Value = '  -0.55%  '
# Litecoin
$ws.Range('D17').Value = '''60.67'
$ws.Range('E17').Value = '  -1.77%  '
# ShibaInu
$ws.Range('D18').Value = '0.0₃0715'
$ws.Range('E18').Value = '  -3.97%  '
# Dai
$ws.Range('E19').Value = '  +0.24%  '
# BitcoinCash
$ws.Range('D20').Value = '''189.03'
$ws.Range('E20').Value = '  -0.71%  '
# Uniswap
$ws.Range('E21').Value = '  -1.85%  '
# Avalanche
$ws.Range('D22').Value = '''9.31'
$ws.Range('E22').Value = '  -2.67%  '
# Chainlink
$ws.Range('D23').Value = '''5.94'
$ws.Range('E23').Value = '  -3.07%  '
# BinanceUSD
$ws.Range('E24').Value = '  +0.02%  '
# Monero
$ws.Range('D25').Value = '''141.72'
$ws.Range('E25').Value = '  -1.16%  '
# Stellar
$ws.Range('E26').Value = '  -3.56%  '
# Toncoin
$ws.Range('E27').Value = '  -3.30%  '
# Cosmos
$ws.Range('D28').Value = '''6.51'
$ws.Range('E28').Value = '  -3.86%  '
# EthereumClassic
$ws.Range('D29').Value = '''14.91'
$ws.Range('E29').Value = '  -1.92%  '
# PancakeSwap
$ws.Range('E30').Value = '  -2.59%  '
# Hedera
$ws.Range('E31').Value = '  -3.36%  '
# Filecoin
$ws.Range('E32').Value = '  -2.53%  '
# InternetComputer(DFINITY)
$ws.Range('D33').Value = '''3.01'
$ws.Range('E33').Value = '  -4.46%  '
# HuobiToken
$ws.Range('E34').Value = '  -1.49%  '
# LidoDAOToken
$ws.Range('E35').Value = '  -2.78%  '
# Maker
$ws.Range('D36').Value = '1.106.02'
$ws.Range('E36').Value = '  -2.53%  '
# MXToken
$ws.Range('D37').Value = '''2.36'
$ws.Range('E37').Value = '  -2.83%  '
# ARBITRUM
$ws.Range('D38').Value = '''0.797'
$ws.Range('E38').Value = '  -8.32%  '
# VeChain
$ws.Range('D39').Value = '''0.0151'
$ws.Range('E39').Value = '  -2.44%  '
# ImmutableX
$ws.Range('D40').Value = '''0.495'
$ws.Range('E40').Value = '  -5.63%  '
# Quant
$ws.Range('D41').Value = '''95.50'
$ws.Range('E41').Value = '  -3.20%  '
# RocketPoolETH
$ws.Range('D42').Value = '1.733.19'
$ws.Range('E42').Value = '  -2.10%  '
# FraxShare
$ws.Range('E43').Value = '  -3.59%  '
# TrustWalletToken
$ws.Range('D44').Value = '''0.739'
$ws.Range('E44').Value = '  -5.33%  '
# BabyDogeCoin
$ws.Range('E45').Value = '  -1.03%  '
# Aave
$ws.Range('D46').Value = '''52.97'
$ws.Range('E46').Value = '  -3.88%  '
# row47 (was Cronos -> now RenderToken)
$ws.Range('B47').Value = 'RenderToken'
$ws.Range('C47').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D47').Value = '''1.46'
$ws.Range('E47').Value = '  -1.17%  '
# row48 (was RenderToken -> now Cronos)
$ws.Range('B48').Value = 'Cronos'
$ws.Range('C48').Value = 'https://coinranking.com/coin/65PHZTpmE55b+cronos-cro'
$ws.Range('D48').Value = '''0.0509'
$ws.Range('E48').Value = '  -3.31%  '
# Mantle
$ws.Range('E49').Value = '  -1.10%  '
# USDD
$ws.Range('E50').Value = '  +0.11%  '
# EnergySwap
$ws.Range('D51').Value = '''7.32'
$ws.Range('E51').Value = '  -2.80%  '
